$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 23
$ws.Range("F7").Value = 1936
$ws.Range("F8").Value = 5420
$ws.Range("F9").Value = 1523
$ws.Range("F10").Value = 155
$ws.Range("F11").Value = 3108
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 1291
$ws.Range("F15").Value = 4290
$ws.Range("F16").Value = 1029
$ws.Range("F18").Value = 1672
$ws.Range("F21").Value = 29
$ws.Range("F24").Value = 976
$ws.Range("F25").Value = 291
$ws.Range("F27").Value = 82
$ws.Range("F29").Value = 1088
$ws.Range("F30").Value = 385
$ws.Range("F31").Value = 51
$ws.Range("F32").Value = 159
$ws.Range("F34").Value = 287
$ws.Range("F36").Value = 1673
$ws.Range("F38").Value = 1019
$ws.Range("F40").Value = 250
$ws.Range("F41").Value = 609
$ws.Range("F42").Value = 293
$ws.Range("C43").Value = "杭州·ACG-World×梦漫星河动漫嘉年华"
$ws.Range("D43").Value = "杭州滨江长江南路336号 白马湖国际会展中心"
$ws.Range("E43").Value = "2024.07.20 10:00-07.21 17:00"
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 68
$ws.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=84813"
$ws.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202404/OVTgtrwI1713457945698.jpeg"
$ws.Range("F44").Value = 655
$ws.Range("F45").Value = 6
$ws.Range("F47").Value = 336
$ws.Range("F49").Value = 139
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 34
$ws.Range("F10").Value = 148
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 744
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 744
$ws.Range("F4").Value = 23
$ws.Range("F8").Value = 1936
$ws.Range("F9").Value = 5420
$ws.Range("F10").Value = 1523
$ws.Range("F11").Value = 155
$ws.Range("F13").Value = 3109
$ws.Range("F14").Value = 40
$ws.Range("F15").Value = 1291
$ws.Range("F16").Value = 4290
$ws.Range("F17").Value = 1029
$ws.Range("F18").Value = 1672
$ws.Range("F20").Value = 34
$ws.Range("F24").Value = 29
$ws.Range("F26").Value = 148
$ws.Range("F27").Value = 976
$ws.Range("F28").Value = 291
$ws.Range("F29").Value = 82
$ws.Range("F31").Value = 1088
$ws.Range("F32").Value = 385
$ws.Range("F33").Value = 51
$ws.Range("F34").Value = 159
$ws.Range("F36").Value = 1673
$ws.Range("F38").Value = 1019
$ws.Range("F42").Value = 250
$ws.Range("F43").Value = 609
$ws.Range("F44").Value = 293
$ws.Range("F45").Value = 655
$ws.Range("F47").Value = 336
$ws.Range("F49").Value = 139
Write-Host "Edits applied successfully"
